# Auto-generated edit script applying the cryptos.xlsx price/volume refresh.
# Forces each touched cell to Text format before assignment so that
# numeric-looking strings (e.g. "1.00", "27.028.08") are preserved
# verbatim as text rather than being coerced to numbers by Excel's
# automatic type detection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.028.08'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.561.40'

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.46%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.23%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.34'

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.69%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.489'

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.51%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.38%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.05'

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.37%  '

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.82%  '

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.76%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.33%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.785.62'

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.51%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.569.61'

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.92%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.73'

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.38%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.520'

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.19%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.030.58'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.88'

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.32%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0707'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.40%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '215.61'

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.05%  '

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.10%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.29%  '

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.77%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.20'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.28%  '

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.51%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.95'

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.98%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.60'

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.47%  '

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.69%  '

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.45%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.22%  '

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.19%  '

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +3.29%  '

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.04%  '

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.53%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.429.52'

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.25%  '

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.66%  '

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +9.20%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.33'

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.32%  '

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.89%  '

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.14%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.89'

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.08%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.809'

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.24%  '

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.36%  '

$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'WEMIXToken'

$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.99%  '

$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'MXToken'

$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.31'

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.29%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.69'

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.61%  '

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.02%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.698.96'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.48%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '86.91'

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.92%  '

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.54%  '

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.58%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0960'

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.50%  '
